$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new values look numeric need to be forced to
# Text format first, otherwise Excel auto-converts them to numbers (the
# source data keeps them as literal text, matching the original file).
$textForceCells = @("D5","D6","D10","D11","D20","D21","D25","D27","D31","D32","D35","D36","D37","D38","D40","D41","D42","D44","D48")
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '67.967.38'
$ws.Range("E2").Value = '  +3.35%  '
$ws.Range("D3").Value = '3.280.87'
$ws.Range("E3").Value = '  +3.42%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '582.55'
$ws.Range("E5").Value = '  +1.69%  '
$ws.Range("D6").Value = '182.75'
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("E8").Value = '  +1.26%  '
$ws.Range("D9").Value = '3.278.60'
$ws.Range("E9").Value = '  +3.45%  '
$ws.Range("D10").Value = '0.135'
$ws.Range("E10").Value = '  +7.51%  '
$ws.Range("D11").Value = '6.73'
$ws.Range("E11").Value = '  +1.86%  '
$ws.Range("E12").Value = '  +6.37%  '
$ws.Range("D13").Value = '3.846.83'
$ws.Range("E13").Value = '  +3.38%  '
$ws.Range("E14").Value = '  +1.37%  '
$ws.Range("E15").Value = '  +3.80%  '
$ws.Range("D16").Value = '67.941.08'
$ws.Range("E16").Value = '  +3.34%  '
$ws.Range("E17").Value = '  +3.20%  '
$ws.Range("D18").Value = '3.281.87'
$ws.Range("E18").Value = '  +3.35%  '
$ws.Range("E19").Value = '  +1.84%  '
$ws.Range("D20").Value = '13.55'
$ws.Range("E20").Value = '  +4.52%  '
$ws.Range("D21").Value = '377.24'
$ws.Range("E21").Value = '  +4.33%  '
$ws.Range("E22").Value = '  +5.37%  '
$ws.Range("E23").Value = '  -0.04%  '
$ws.Range("E24").Value = '  +3.09%  '
$ws.Range("D25").Value = '0.514'
$ws.Range("E25").Value = '  +3.74%  '
$ws.Range("E26").Value = '  +5.59%  '
$ws.Range("D27").Value = '9.72'
$ws.Range("E27").Value = '  -1.31%  '
$ws.Range("E28").Value = '  +2.31%  '
$ws.Range("E29").Value = '  +0.06%  '
$ws.Range("E30").Value = '  +2.82%  '
$ws.Range("D31").Value = '5.71'
$ws.Range("E31").Value = '  +5.07%  '
$ws.Range("D32").Value = '22.92'
$ws.Range("E32").Value = '  +3.55%  '
$ws.Range("E33").Value = '  +6.78%  '
$ws.Range("D35").Value = '6.95'
$ws.Range("E35").Value = '  +4.52%  '
$ws.Range("D36").Value = '1.54'
$ws.Range("E36").Value = '  +5.22%  '
$ws.Range("D37").Value = '162.32'
$ws.Range("E37").Value = '  +1.20%  '
$ws.Range("D38").Value = '0.854'
$ws.Range("E38").Value = '  +2.04%  '
$ws.Range("E39").Value = '  +2.84%  '
$ws.Range("D40").Value = '27.08'
$ws.Range("E40").Value = '  +2.56%  '
$ws.Range("D41").Value = '6.79'
$ws.Range("E41").Value = '  +9.89%  '
$ws.Range("D42").Value = '4.62'
$ws.Range("E42").Value = '  +9.92%  '
$ws.Range("E43").Value = '  +5.47%  '
$ws.Range("D44").Value = '352.09'
$ws.Range("E44").Value = '  +6.30%  '
$ws.Range("E45").Value = '  +7.10%  '
$ws.Range("D46").Value = '2.679.98'
$ws.Range("E46").Value = '  +1.08%  '
$ws.Range("E47").Value = '  +2.77%  '
$ws.Range("D48").Value = '0.0682'
$ws.Range("E49").Value = '  +3.19%  '
$ws.Range("E50").Value = '  +5.14%  '
$ws.Range("E51").Value = '  +1.35%  '
